$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Reorganize the "Estado de Cuenta" detail rows (16-29): previously grouped by
# employee (all periods for MIGUEL, then all periods for CARLOS, each in
# descending period order). Now interleaved by period (ascending), alternating
# MIGUEL / CARLOS per period, and the "Salario Basico" for MIGUEL's rows is
# updated from 800000 to 1160000.

$rows = @(
    @{ Row = 16; Doc = "15326489"; Name = "MIGUEL AVILA MARMOL";             Period = "2306"; Mora = 46400; Salario = 1160000 },
    @{ Row = 17; Doc = "80019630"; Name = "CARLOS URIEL CRISTANCHO SALAZAR"; Period = "2306"; Mora = 48000; Salario = 1200000 },
    @{ Row = 18; Doc = "15326489"; Name = "MIGUEL AVILA MARMOL";             Period = "2307"; Mora = 46400; Salario = 1160000 },
    @{ Row = 19; Doc = "80019630"; Name = "CARLOS URIEL CRISTANCHO SALAZAR"; Period = "2307"; Mora = 48000; Salario = 1200000 },
    @{ Row = 20; Doc = "15326489"; Name = "MIGUEL AVILA MARMOL";             Period = "2308"; Mora = 46400; Salario = 1160000 },
    @{ Row = 21; Doc = "80019630"; Name = "CARLOS URIEL CRISTANCHO SALAZAR"; Period = "2308"; Mora = 48000; Salario = 1200000 },
    @{ Row = 22; Doc = "15326489"; Name = "MIGUEL AVILA MARMOL";             Period = "2309"; Mora = 46400; Salario = 1160000 },
    @{ Row = 23; Doc = "80019630"; Name = "CARLOS URIEL CRISTANCHO SALAZAR"; Period = "2309"; Mora = 48000; Salario = 1200000 },
    @{ Row = 24; Doc = "15326489"; Name = "MIGUEL AVILA MARMOL";             Period = "2310"; Mora = 46400; Salario = 1160000 },
    @{ Row = 25; Doc = "80019630"; Name = "CARLOS URIEL CRISTANCHO SALAZAR"; Period = "2310"; Mora = 48000; Salario = 1200000 },
    @{ Row = 26; Doc = "15326489"; Name = "MIGUEL AVILA MARMOL";             Period = "2311"; Mora = 46400; Salario = 1160000 },
    @{ Row = 27; Doc = "80019630"; Name = "CARLOS URIEL CRISTANCHO SALAZAR"; Period = "2311"; Mora = 48000; Salario = 1200000 },
    @{ Row = 28; Doc = "15326489"; Name = "MIGUEL AVILA MARMOL";             Period = "2312"; Mora = 29387; Salario = 1160000 },
    @{ Row = 29; Doc = "80019630"; Name = "CARLOS URIEL CRISTANCHO SALAZAR"; Period = "2312"; Mora = 30400; Salario = 1200000 }
)

foreach ($r in $rows) {
    $i = $r.Row
    $ws.Range("B$i").Value = "CC"
    $ws.Range("C$i").Value = $r.Doc
    $ws.Range("D$i").Value = $r.Name
    $ws.Range("E$i").Value = $r.Period
    $ws.Range("F$i").Value = $r.Mora
    $ws.Range("G$i").Value = $r.Salario
}
